$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column B (existing B,C shift right to C,D)
$ws.Columns.Item(2).Insert()

# Adjust column widths: A grows, new B takes old A's old width, C/D keep old B/C widths
$ws.Columns.Item(1).ColumnWidth = 124.6640625
$ws.Columns.Item(2).ColumnWidth = 85.6640625
$ws.Columns.Item(3).ColumnWidth = 20.5
$ws.Columns.Item(4).ColumnWidth = 27.83203125

# Row 4 & 5: new data rows, set in the order that reproduces the author's
# shared-string table layout (A4, D4, C4, A5, D5, C5, B4, B5)
$ws.Range("A4").Value = "Mapping legend search macropus"
$ws.Range("D4").Value = "q=macropus&cm=species&type=application/json"
$ws.Range("C4").Value = "mapping/legend"

$ws.Range("A5").Value = "Generate static map"
$ws.Range("D5").Value = "pcolour=3531FF&popacity=1&outlineColour=0x000000&dpi=300&scale=on&format=jpg&outline=true&q=Macropus+rufus&extents=96.173828125,-47.11468820158343,169.826171875,-2.5694811631203973&baselayer=world&fileName=MyMap.jpg&pradiusmm=1"
$ws.Range("C5").Value = "mapping/wms/image"

$ws.Range("B4").Value = "application/json"
$ws.Range("B5").Value = "application/octet-stream"

# Row 5, columns C & D get a 10pt Arial font (this style is registered
# first, becoming cellXfs index 1). Build it once on C5, then copy the
# format across to D5 so both share the same style entry.
$ws.Range("C5").Font.Size = 10
$ws.Range("C5").Font.Name = "Arial"
$ws.Range("C5").Copy()
$ws.Range("D5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New column B (rows 1-3) is blank but carries a small Monaco/green font
# (registered second, becoming cellXfs index 2). Build the style once on
# B1, then fan it out via copy/paste-format so every cell shares a single
# style entry instead of re-deriving the font per cell.
$ws.Range("B1").Font.Size = 9
$ws.Range("B1").Font.Color = 5867370
$ws.Range("B1").Font.Name = "Monaco"
$ws.Range("B1").Copy()
$ws.Range("B2:B3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update selection to B5
$ws.Range("B5").Select()

# Best-effort: sync the window view metrics (may not persist in this runtime)
$excel.ActiveWindow.Left = 0
$excel.ActiveWindow.Top = 20
$excel.ActiveWindow.Width = 38400
$excel.ActiveWindow.Height = 21600
